$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.519.33"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.949.99"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.82"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.622"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.54"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.80%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.367"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.46%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0842"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.89%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.80"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.91%  "
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.827"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.232.63"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.65"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.947.22"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.393.17"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.81"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0865"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.24"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.06"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.74%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.45"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.30"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.24"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -5.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.99"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.135"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.52"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.29%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.16"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.70"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0632"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.29"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.29"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.21%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.78"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.14"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -5.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.07"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0987"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.18"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0210"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.05"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.367.59"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.04"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.24"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.16"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.81%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.85"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.05%  "
